# Update column F ("dSF") values for the data rows.
# Row 7 is intentionally left unchanged (stays at 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = -4
    4  = -3
    5  = 1
    6  = -1
    8  = 1
    9  = 2
    10 = 1
    11 = -3
    12 = 1
    13 = -6
    14 = 2
    15 = -3
    16 = 1
    17 = 9
    18 = 2
    19 = 3
    20 = 0
    21 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
